$d = $word.ActiveDocument

$names = @(
    "FLX3D518BAF",
    "FLXF178377",
    "FLX520AEBFC",
    "FLX4565B83C",
    "FLX344ACD8A",
    "FLX7ACDF45D",
    "FLX13B390B0",
    "FLX2314008D",
    "FLX465323BA",
    "FLX658FCE2D",
    "FLX1B44E5FC",
    "FLX63872413",
    "FLX1F6F159C",
    "FLX37C67B77"
)

$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $para = $d.Paragraphs($i)
    $r = $para.Range
    $r.End = $r.End - 1
    $name = $names[$i - 1]
    $d.Bookmarks.Add($name, $r)
}

Write-Output "done"
